$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 7 (week 4) topic: Data-Wrangling -> Data-Visualization
$ws.Range("C7").Value = "[Data-Visualization](https://crumplab.github.io/psyc7709/Schedule.html#4_data-visualization)"

# Row 8 (week 5) topic: Data-Visualization -> Data-Wrangling
$ws.Range("C8").Value = "[Data-Wrangling](https://crumplab.github.io/psyc7709/Schedule.html#5_data-wrangling)"

# Row 8 (week 5) assignment due label: data-wrangling -> data-vis
$ws.Range("D8").Value = "week 4 due (data-vis)"

# Row 9 (week 6) assignment due label: data-vis -> data-wrangling
$ws.Range("D9").Value = "week 5 due (data-wrangling)"

# Update selection / scroll position to match the saved view state.
$ws.Range("D10").Select()
